# "Generate Report for Archive"
#
# Refresh the localization-status report: every item that was previously
# marked "Ready for handoff" has since moved on to "In Translation", so
# update the status wherever it is shown (the per-language status columns
# on the Overview sheet, and the "Status" column on each language sheet),
# then re-tighten those columns now that the status text is shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

function Update-StatusColumn {
    param($ws, $colIndex, $firstRow, $lastRow)

    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $colIndex)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Overview sheet ---------------------------------------------------
# Column E = "zh-cn" status, Column F = "de-de" status; data rows 2-4.
$overview = $wb.Worksheets.Item("Overview")
Update-StatusColumn $overview 5 2 4
Update-StatusColumn $overview 6 2 4
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-language detail sheets ---------------------------------------
# Column C = "Status"; data rows 2-4 on both the zh-cn and de-de sheets.
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    Update-StatusColumn $ws 3 2 4
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
